# Solve Leetcode - 286. Walls and Gates - Multi Source BFS and DFS
# Adds two new rows to the "Neetcode 150" journal sheet for:
#   133. Clone Graph
#   286. Walls and Gates

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new note text -------------------------------------------------------

$cloneGraphTitle = "133. Clone Graph"
$cloneGraphNote = @'
We need to create a copy for each node in the graph, and we'll track those copies with a hashap
Then, run a dfs with base case node is already in hashmap then just return it from the hashmap otherwise create the node, add it to hashmap and run dfs again recursively for each neighbour and append what it returns to copy's neighbors, and return the copy after going through all neighbors
'@

$wallsTitle = "286. Walls and Gates"
$wallsNote = @'
Multi Source BFS is just like bfs except we start with multiple root nodes and move deeper for all nodes consecutively.
First find all 0's, those will be the roots for the multi source bfs, then run the usual bfs exapnding in all 4 directions, and for each direction, add it to bfs if its indices are within bounds and if its value is <= current val + 1, only then overwrite its dist and add it to null to ensure that we only overwrite to make paths shorter, not longer.
'@

# ---- row 32: 133. Clone Graph --------------------------------------------

$ws.Range("A32").Value = "Graphs"
$ws.Range("B32").Value = "Medium"

# Add the hyperlink first (TextToDisplay mirrors the URL, matching the way
# the other rows in this sheet were authored), then overwrite the cell's
# visible text with the real problem title, and finally restore the normal
# "Good" cell formatting that Hyperlinks.Add clobbers.
$ws.Hyperlinks.Add($ws.Range("C32"), "https://leetcode.com/problems/clone-graph/", "", "", "https://leetcode.com/problems/clone-graph/")
$ws.Range("C32").Value = $cloneGraphTitle
$ws.Range("C32").Style = "Good"

$ws.Range("D32").Value = $cloneGraphNote
$ws.Range("D32").WrapText = $true
$ws.Range("D32").VerticalAlignment = -4160

$ws.Range("B32").Style = "Neutral"

$ws.Rows.Item(32).RowHeight = 57.6

# ---- row 33: 286. Walls and Gates ----------------------------------------

$ws.Range("A33").Value = "Graphs"
$ws.Range("B33").Value = "Medium"

$ws.Hyperlinks.Add($ws.Range("C33"), "https://leetcode.com/problems/walls-and-gates/")
$ws.Range("C33").Value = $wallsTitle
$ws.Range("C33").Style = "Good"

$ws.Range("D33").Value = $wallsNote
$ws.Range("D33").WrapText = $true
$ws.Range("D33").VerticalAlignment = -4160

$ws.Range("B33").Style = "Neutral"

$ws.Rows.Item(33).RowHeight = 57.6

# ---- selection / scroll position ------------------------------------------

[void]$ws.Range("A34").Select()
